# Weekly update: insert a new price-report row for "Feria Lagunitas de
# Puerto Montt" (Ajo / Chino) at row 81, shifting the existing rows
# 81-156 down to 82-157.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 81, shifting rows 81..156 down to
# 82..157 (xlShiftDown = -4121).
$ws.Rows(81).Insert(-4121)

# Populate the newly inserted row with the latest week's data.
$ws.Cells.Item(81, 1).Value  = 4
$ws.Cells.Item(81, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(81, 3).Value  = "Los Lagos"
$ws.Cells.Item(81, 4).Value  = 44484
$ws.Cells.Item(81, 5).Value  = 10
$ws.Cells.Item(81, 6).Value  = 100112003
$ws.Cells.Item(81, 7).Value  = "Ajo"
$ws.Cells.Item(81, 8).Value  = "Chino"
$ws.Cells.Item(81, 9).Value  = "Primera"
$ws.Cells.Item(81, 10).Value = 240
$ws.Cells.Item(81, 11).Value = 17500
$ws.Cells.Item(81, 12).Value = 19000
$ws.Cells.Item(81, 13).Value = 18250
$ws.Cells.Item(81, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(81, 15).Value = "China"
$ws.Cells.Item(81, 16).Value = 1825
$ws.Cells.Item(81, 17).Value = 10
$ws.Cells.Item(81, 18).Value = "Hortaliza"
